$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("questions")

# Update the "Others" module time question row (row 8):
# - B8: "Pick a time (24 hrs)" -> "Pick a time"
# - D8: (empty) -> "09:00,17:00,12:00"
$ws.Range("B8").Value = "Pick a time"
$ws.Range("D8").Value = "09:00,17:00,12:00"

# Update the selected/active cell to D8, matching the post-edit UI state
$ws.Range("D8").Select()
